$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update generated timestamp string in A2
$ws.Range("A2").Value = "Előállítva: 11/24/2025 6:00:19 AM, Készítette: Gazso"

# Update quantity values in column B
$ws.Range("B120").Value = 5
$ws.Range("B121").Value = 5
$ws.Range("B122").Value = 5
$ws.Range("B136").Value = 1
$ws.Range("B137").Value = 1
$ws.Range("B138").Value = 1
$ws.Range("B139").Value = 1
$ws.Range("B140").Value = 1
$ws.Range("B141").Value = 1
$ws.Range("B175").Value = 38
$ws.Range("B176").Value = 38
$ws.Range("B177").Value = 38
$ws.Range("B209").Value = 6
$ws.Range("B210").Value = 6
$ws.Range("B211").Value = 6
$ws.Range("B215").Value = 11
$ws.Range("B216").Value = 11
$ws.Range("B217").Value = 11
$ws.Range("B218").Value = 4
$ws.Range("B219").Value = 4
$ws.Range("B220").Value = 4
$ws.Range("B221").Value = 18
$ws.Range("B222").Value = 18
$ws.Range("B223").Value = 0
$ws.Range("B279").Value = 8
$ws.Range("B280").Value = 8
$ws.Range("B281").Value = 8
$ws.Range("B324").Value = 0
$ws.Range("B325").Value = 0
$ws.Range("B326").Value = 0
$ws.Range("B358").Value = 2
$ws.Range("B359").Value = 2
$ws.Range("B360").Value = 2
$ws.Range("B367").Value = 3
$ws.Range("B368").Value = 3
$ws.Range("B369").Value = 3
$ws.Range("B626").Value = 1
$ws.Range("B627").Value = 1
$ws.Range("B628").Value = 1
$ws.Range("B647").Value = 30
$ws.Range("B648").Value = 30
$ws.Range("B649").Value = 30
$ws.Range("B666").Value = 4
$ws.Range("B667").Value = 4
$ws.Range("B668").Value = 4
$ws.Range("B672").Value = 13
$ws.Range("B673").Value = 13
$ws.Range("B674").Value = 13
$ws.Range("B675").Value = 6
$ws.Range("B676").Value = 6
$ws.Range("B677").Value = 6
$ws.Range("B715").Value = 4
$ws.Range("B716").Value = 4
$ws.Range("B717").Value = 4
$ws.Range("B738").Value = 3
$ws.Range("B739").Value = 3
$ws.Range("B740").Value = 3
$ws.Range("B773").Value = 9
$ws.Range("B774").Value = 9
$ws.Range("B775").Value = 9
$ws.Range("B801").Value = 1
$ws.Range("B802").Value = 1
$ws.Range("B803").Value = 1
$ws.Range("B807").Value = 8
$ws.Range("B808").Value = 8
$ws.Range("B809").Value = 8
$ws.Range("B813").Value = 1
$ws.Range("B814").Value = 1
$ws.Range("B815").Value = 1
$ws.Range("B825").Value = 2
$ws.Range("B826").Value = 2
$ws.Range("B827").Value = 2

Write-Host "Applied 72 cell updates"
